$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header info
$ws.Range("C2").Value = "Hartmut"
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 10.06.2025"

# Row 6 (transaction 1)
$ws.Range("B6").Value = "13.06."
$ws.Range("C6").Value = "14.06."
$ws.Range("D6").Value = "ZALANDO MKTPLC EU FDDCBB"
$ws.Range("E6").Value = "203,34-"

# Row 7 (transaction 2)
$ws.Range("B7").Value = "16.06."
$ws.Range("C7").Value = "17.06."
$ws.Range("D7").Value = "RECHNUNG VODAFONE GMBH 18324118"
$ws.Range("E7").Value = "39,80-"

# Row 8 (transaction 3)
$ws.Range("B8").Value = "17.06."
$ws.Range("C8").Value = "18.06."
$ws.Range("D8").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 13899023"
$ws.Range("E8").Value = "83,95-"

# Row 9 (transaction 4) removed entirely - clear contents; the amount cell
# E9 also switches from right-aligned to center-aligned (matches the blank
# trailing rows 10/11 look) once it no longer holds a value.
$ws.Range("B9").Value = $null
$ws.Range("C9").Value = $null
$ws.Range("D9").Value = $null
$ws.Range("E9").Value = $null
$ws.Range("E9").HorizontalAlignment = -4108  # xlCenter
$ws.Range("E9").VerticalAlignment = -4108    # xlCenter
$ws.Range("E9").WrapText = $true

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 20.06.2025"
$ws.Range("E12").Value = "327,09-"

# Next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 28.06.2025"
